# Updated cryptos list on Tue Jul  4 14:30:47 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "31.011.48"
$ws.Range("E2").Value = "  +1.14%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.958.56"
$ws.Range("E3").Value = "  -0.11%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.67"
$ws.Range("E5").Value = "  -1.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9996"
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4868"
$ws.Range("E7").Value = "  +1.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2952"
$ws.Range("E8").Value = "  +0.85%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06937"
$ws.Range("E9").Value = "  +2.66%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.50"
$ws.Range("E10").Value = "  +1.83%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "108.37"
$ws.Range("E11").Value = "  -0.60%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.965.88"
$ws.Range("E12").Value = "  +0.25%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07806"
$ws.Range("E13").Value = "  +0.88%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.504"
$ws.Range("E14").Value = "  +0.70%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.7027"
$ws.Range("E15").Value = "  +1.14%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "281.78"
$ws.Range("E16").Value = "  -3.59%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "31.025.90"
$ws.Range("E17").Value = "  +1.13%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.33"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007786"
$ws.Range("E19").Value = "  +1.31%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.236.27"
$ws.Range("E20").Value = "  +0.92%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9992"
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.540"
$ws.Range("E22").Value = "  -2.42%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9998"
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.538"
$ws.Range("E24").Value = "  -1.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.889"
$ws.Range("E25").Value = "  -0.25%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "168.76"
$ws.Range("E26").Value = "  -0.41%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.99"
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.198"
$ws.Range("E29").Value = "  -1.55%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.386"
$ws.Range("E30").Value = "  -3.89%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.665"
$ws.Range("E31").Value = "  -4.10%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.570"
$ws.Range("E32").Value = "  -1.92%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.483"
$ws.Range("E33").Value = "  +0.47%  "
$ws.Range("E34").Value = "  -3.09%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7589"
$ws.Range("E35").Value = "  -1.64%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.173"
$ws.Range("E36").Value = "  -0.45%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.729"
$ws.Range("E37").Value = "  +0.24%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02020"
$ws.Range("E38").Value = "  -0.89%  "
$ws.Range("E39").Value = "  -0.47%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.600"
$ws.Range("E40").Value = "  +1.79%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "78.16"
$ws.Range("E41").Value = "  +11.56%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.129"
$ws.Range("E42").Value = "  -0.21%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9031"
$ws.Range("E43").Value = "  +2.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "109.84"
$ws.Range("E44").Value = "  -0.16%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4465"
$ws.Range("E45").Value = "  +0.00%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.167"
$ws.Range("E46").Value = "  +9.28%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.9998"
$ws.Range("E47").Value = "  +0.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.013.57"
$ws.Range("E48").Value = "  +8.73%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.376"
$ws.Range("E49").Value = "  -0.12%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1259"
$ws.Range("E50").Value = "  -1.47%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "36.04"
$ws.Range("E51").Value = "  +0.05%  "
